$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new annotation row (row 6) for "paris" with the same layout as
# the existing rows: Annotator, politeness_score, polite_expressions,
# sentence_purpose, issue_type, id, source_file, text
$ws.Cells.Item(6, 1).Value = "paris"
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = "CRT"
$ws.Cells.Item(6, 5).Value = "RES"
$ws.Cells.Item(6, 6).Value = "dc9804e9-fe90-49ab-88bb-ac97478c1b97"
$ws.Cells.Item(6, 7).Value = "i87JIQTAnB8AQ_annotated.xlsx"
$ws.Cells.Item(6, 8).Value = "As you suggested, I did run comparison tests and I will present the results here."
